$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.159.08'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.51%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.503.56'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.88%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.55'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.82%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '94.33'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.06%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.568'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.95%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.14%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.526'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.40%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.86'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -3.06%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0800'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.66%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.35'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -2.24%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -5.20%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.897.24'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.54%  '
$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.527.08'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.12%  '
$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.79'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -5.01%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.829'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -3.17%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.443.82'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.08%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.71'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.97%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.29'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -4.72%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0945'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -1.64%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.72'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -2.89%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '249.79'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.11%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.91'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -1.61%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.01'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.83%  '
$ws.Range("B26").Value = 'Dai'
$ws.Range("C26").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.42%  '
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '26.23'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -2.25%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.39'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.66%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '39.56'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.48%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '10.07'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.05%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.75'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -2.91%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '155.12'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.08%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.69'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +3.03%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '18.95'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +1.17%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.25'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -1.44%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0780'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.11%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.05'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -3.79%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.110'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.33%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.43'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +2.65%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -1.47%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '21.84'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -8.11%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.18%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0298'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.05%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.72'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -2.35%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.19'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -3.99%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.975.11'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -2.55%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.92'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +1.11%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '83.03'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -1.48%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.766.03'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.04%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '102.63'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +1.34%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.84'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.34%  '
